$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 669.5
$ws.Range("J45").Value = 669.5
$ws.Range("L45").Value = 2008.5
$ws.Range("N45").Value = -2392.5

$ws.Range("H46").Value = 300
$ws.Range("J46").Value = 300
$ws.Range("L46").Value = 900
$ws.Range("N46").Value = -1138

$ws.Range("H60").Value = 300
$ws.Range("J60").Value = 300
$ws.Range("L60").Value = 900
$ws.Range("N60").Value = -1868

$ws.Range("H62").Value = 6946989
$ws.Range("I62").Value = 9261326
$ws.Range("K62").Value = 9261326
$ws.Range("M62").Value = -9260702

$ws.Range("H64").Value = 3750.6
$ws.Range("J64").Value = 3440.75
$ws.Range("L64").Value = 3440.75
$ws.Range("N64").Value = -3936.75

$ws.Range("H65").Value = 6946989
$ws.Range("I65").Value = 9261326
$ws.Range("K65").Value = 46306630
$ws.Range("M65").Value = -46303510

$ws.Range("H67").Value = 3750.6
$ws.Range("J67").Value = 3440.75
$ws.Range("L67").Value = 3440.75
$ws.Range("N67").Value = -5156.75

$ws.Range("H76").Value = 5370.857
$ws.Range("I76").Value = 4199.3335
$ws.Range("K76").Value = 4199.3335
$ws.Range("M76").Value = -3884.3335

$ws.Range("H79").Value = 5370.857
$ws.Range("I79").Value = 4199.3335
$ws.Range("K79").Value = 4199.3335
$ws.Range("M79").Value = -3107.3335

$ws.Range("H116").Value = 2603.4285
$ws.Range("I116").Value = 2356.125
$ws.Range("J116").Value = 2933.1667
$ws.Range("K116").Value = 2356.125
$ws.Range("L116").Value = 2933.1667
$ws.Range("M116").Value = 1085.875
$ws.Range("N116").Value = -9817.1667

$ws.Range("H134").Value = 36460
$ws.Range("J134").Value = 36460
$ws.Range("L134").Value = 36460
$ws.Range("N134").Value = -46600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 323.5
$ws.Range("I5").Value = 281.33334
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 281.33334
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = -169.33334
$ws.Range("N5").Value = -674

$ws.Range("H45").Value = 971.0952
$ws.Range("I45").Value = 829
$ws.Range("K45").Value = 829
$ws.Range("M45").Value = -452

$ws.Range("H63").Value = 2433.8965
$ws.Range("I63").Value = 2392.5186
$ws.Range("J63").Value = 2992.5
$ws.Range("K63").Value = 2392.5186
$ws.Range("L63").Value = 2992.5
$ws.Range("M63").Value = -1706.5186
$ws.Range("N63").Value = -4364.5

$ws.Range("H66").Value = 2433.8965
$ws.Range("I66").Value = 2392.5186
$ws.Range("J66").Value = 2992.5
$ws.Range("K66").Value = 11962.593
$ws.Range("L66").Value = 14962.5
$ws.Range("M66").Value = -8530.592999999999
$ws.Range("N66").Value = -21826.5

$ws.Range("H74").Value = 2997.7778
$ws.Range("I74").Value = 2311.4285
$ws.Range("K74").Value = 2311.4285
$ws.Range("M74").Value = -1437.4285

$ws.Range("H77").Value = 2997.7778
$ws.Range("I77").Value = 2311.4285
$ws.Range("K77").Value = 11557.1425
$ws.Range("M77").Value = -7189.1425

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 323.5
$ws.Range("I4").Value = 281.33334
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 281.33334
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -166.33334
$ws.Range("N4").Value = -680

$ws.Range("H82").Value = 15645.2
$ws.Range("I82").Value = 2326.75
$ws.Range("J82").Value = 30866.285
$ws.Range("K82").Value = 2326.75
$ws.Range("L82").Value = 30866.285
$ws.Range("M82").Value = -1943.75
$ws.Range("N82").Value = -31632.285

$ws.Range("H85").Value = 15645.2
$ws.Range("I85").Value = 2326.75
$ws.Range("J85").Value = 30866.285
$ws.Range("K85").Value = 2326.75
$ws.Range("L85").Value = 30866.285
$ws.Range("M85").Value = -1000.75
$ws.Range("N85").Value = -33518.285

$ws.Range("H105").Value = 250026100
$ws.Range("I105").Value = 250026100
$ws.Range("K105").Value = 250026100
$ws.Range("M105").Value = -250024353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.2
$ws.Range("J2").Value = 66.75
$ws.Range("L2").Value = 400.5
$ws.Range("N2").Value = -626.5

$ws.Range("H7").Value = 496.13333
$ws.Range("I7").Value = 503.30768
$ws.Range("J7").Value = 449.5
$ws.Range("K7").Value = 1509.92304
$ws.Range("L7").Value = 1348.5
$ws.Range("M7").Value = -1397.92304
$ws.Range("N7").Value = -1572.5

$ws.Range("H39").Value = 4338
$ws.Range("J39").Value = 4517.3335
$ws.Range("L39").Value = 13552.0005
$ws.Range("N39").Value = -14140.0005

$ws.Range("H87").Value = 1416.6
$ws.Range("I87").Value = 777.6667
$ws.Range("J87").Value = 2375
$ws.Range("K87").Value = 2333.0001
$ws.Range("L87").Value = 7125
$ws.Range("M87").Value = -1085.0001
$ws.Range("N87").Value = -9621

$ws.Range("H90").Value = 1416.6
$ws.Range("I90").Value = 777.6667
$ws.Range("J90").Value = 2375
$ws.Range("K90").Value = 6999.0003
$ws.Range("L90").Value = 21375
$ws.Range("M90").Value = -759.0002999999997
$ws.Range("N90").Value = -33855

$ws.Range("H131").Value = 17571056
$ws.Range("J131").Value = 36716
$ws.Range("L131").Value = 110148
$ws.Range("N131").Value = -120228

$ws.Range("H138").Value = 3013.5806
$ws.Range("I138").Value = 3144.7
$ws.Range("K138").Value = 9434.099999999999
$ws.Range("M138").Value = -4294.099999999999

$ws.Range("H139").Value = 2173.7407
$ws.Range("I139").Value = 2232.9565
$ws.Range("K139").Value = 6698.869499999999
$ws.Range("M139").Value = -1558.869499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5226
$ws.Range("M6").ClearContents()

$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5500
$ws.Range("M16").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1623.1111
$ws.Range("I16").Value = 1094.4
$ws.Range("J16").Value = 4266.6665
$ws.Range("K16").Value = 1094.4
$ws.Range("L16").Value = 4266.6665
$ws.Range("M16").Value = -924.4000000000001
$ws.Range("N16").Value = -4606.6665

$ws.Range("H93").Value = 1000.3
$ws.Range("I93").Value = 1043.2858
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 1043.2858
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 204.7141999999999
$ws.Range("N93").Value = -3396

$ws.Range("H123").Value = 40482.5
$ws.Range("J123").Value = 40482.5
$ws.Range("L123").Value = 40482.5
$ws.Range("N123").Value = -50282.5

$ws.Range("H132").Value = 2708.75
$ws.Range("I132").Value = 2827.0908
$ws.Range("J132").Value = 2608.6155
$ws.Range("K132").Value = 8481.2724
$ws.Range("L132").Value = 7825.8465
$ws.Range("M132").Value = -5951.2724
$ws.Range("N132").Value = -12885.8465

$ws.Range("H136").Value = 1450.3077
$ws.Range("I136").Value = 783.2222
$ws.Range("K136").Value = 2349.6666
$ws.Range("M136").Value = 200.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1272.8182
$ws.Range("I81").Value = 858.1667
$ws.Range("J81").Value = 1770.4
$ws.Range("K81").Value = 1716.3334
$ws.Range("L81").Value = 3540.8
$ws.Range("M81").Value = -655.3334
$ws.Range("N81").Value = -5662.8

$ws.Range("H84").Value = 1272.8182
$ws.Range("I84").Value = 858.1667
$ws.Range("J84").Value = 1770.4
$ws.Range("K84").Value = 8581.666999999999
$ws.Range("L84").Value = 17704
$ws.Range("M84").Value = -3277.666999999999
$ws.Range("N84").Value = -28312

$ws.Range("H126").Value = 38462420
$ws.Range("I126").Value = 52632292
$ws.Range("J126").Value = 1346.4286
$ws.Range("K126").Value = 157896876
$ws.Range("L126").Value = 4039.2858
$ws.Range("M126").Value = -157894406
$ws.Range("N126").Value = -8979.2858

$ws.Range("H132").Value = 2588.077
$ws.Range("I132").Value = 2114.6
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 6343.799999999999
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -3813.799999999999
$ws.Range("N132").Value = -17559.0005
